# Fix bugs in the GRU example worksheet: the P19/Q19/P20/Q20 formulas were
# referencing the pre-activation row (19/20) instead of the activated
# "reset gate" row (23/24) for the M/N terms.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P19").Formula = "=`$J`$4*G8+`$K`$4*G9+`$L`$4*G10+`$M`$4*G11+M23*(`$A`$27*P8+`$B`$27*P9)+G14"
$ws.Range("Q19").Formula = "=`$J`$4*H8+`$K`$4*H9+`$L`$4*H10+`$M`$4*H11+N23*(`$A`$27*Q8+`$B`$27*Q9)+H14"
$ws.Range("P20").Formula = "=`$J`$5*G8+`$K`$5*G9+`$L`$5*G10+`$M`$5*G11+M24*(`$A`$28*P8+`$B`$28*P9)+G15"
$ws.Range("Q20").Formula = "=`$J`$5*H8+`$K`$5*H9+`$L`$5*H10+`$M`$5*H11+N24*(`$A`$28*Q8+`$B`$28*Q9)+H15"

# Restore the active selection to what the user had selected afterwards.
$ws.Range("Q23").Select()
